# This edit reshuffles the contents of rows 4-15 on the active sheet.
# Every row keeps its row-level formatting / constant columns in place;
# only the columns whose values actually differ row-to-row
# (A, B, D, E, F, G, H, P, Q, R, AO) are moved between rows, following
# the permutation observed between the "before" and "after" workbook
# snapshots. Reading every source value up front (before any writes)
# means the permutation is safe to apply even though some destination
# rows are also source rows for other destinations.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A", "B", "D", "E", "F", "G", "H", "P", "Q", "R", "AO")

# destination row -> source row (source row's data moves into destination row)
$mapping = @{
    4  = 7
    5  = 9
    6  = 13
    7  = 14
    8  = 15
    9  = 4
    10 = 5
    11 = 6
    12 = 8
    13 = 10
    14 = 11
    15 = 12
}

# Snapshot every value we might need to move, before any writes happen.
$snapshot = @{}
foreach ($row in $mapping.Values) {
    if (-not $snapshot.ContainsKey($row)) {
        $rowValues = @{}
        foreach ($c in $cols) {
            $rowValues[$c] = $ws.Range("$c$row").Value2
        }
        $snapshot[$row] = $rowValues
    }
}

# Now write the snapshotted values into their destination rows.
foreach ($destRow in $mapping.Keys) {
    $srcRow = $mapping[$destRow]
    $rowValues = $snapshot[$srcRow]
    foreach ($c in $cols) {
        $ws.Range("$c$destRow").Value2 = $rowValues[$c]
    }
}
